{"js": "// The document contains two \"<id>...</id>\" markers (e.g. \"<id>p051r_1</id>\"\n// and \"<id>p051r_2</id>\"), each currently split across three separate runs:\n//   run1 \"<id>\"      (Courier New, color 7f6000, sz 18)\n//   run2 \"p051r_N\"   (plain, color 000000)\n//   run3 \"</id>\"     (Courier New, color 7f6000, sz 18)\n// The edit collapses each triple into a single run containing the full\n// \"<id>p051r_N</id>\" text, carrying the Courier-New styling of the\n// surrounding \"<id>\"/\"</id>\" runs.\n\nconst body = context.document.body;\n\n// Find every \"<id>...</id>\" span (as a whole, spanning the 3 runs) and\n// rewrite it in one shot. Office.js's Range.insertText(\"Replace\") rewrites\n// the whole matched range as a single new run (adopting the formatting of\n// the range's first character), which is exactly the merge the diff shows.\nconst ids = [\"p051r_1\", \"p051r_2\"];\n\nfor (const id of ids) {\n  const target = `<id>${id}</id>`;\n  const results = body.search(target, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const found of results.items) {\n    found.insertText(target, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# The document contains two \"<id>...</id>\" markers (e.g. \"<id>p051r_1</id>\"\n# and \"<id>p051r_2</id>\"), each currently split across three separate runs:\n#   run1 \"<id>\"      (Courier New, color 7f6000, sz 18)\n#   run2 \"p051r_N\"   (plain, color 000000)\n#   run3 \"</id>\"     (Courier New, color 7f6000, sz 18)\n# The edit collapses each triple into a single run containing the full\n# \"<id>p051r_N</id>\" text, carrying the Courier-New styling of the\n# surrounding \"<id>\"/\"</id>\" runs.\n\n$d = $word.ActiveDocument\n$ids = @(\"p051r_1\", \"p051r_2\")\n\nforeach ($id in $ids) {\n    $target = \"<id>\" + $id + \"</id>\"\n\n    $r = $d.Content\n    $find = $r.Find\n    $find.ClearFormatting()\n    $find.Text = $target\n    $find.MatchCase = $true\n    $find.MatchWildcards = $false\n    $find.Execute() | Out-Null\n\n    if ($find.Found) {\n        # Re-assigning Range.Text to a *different* value first forces Word to\n        # collapse the multi-run span into a single new run that takes on the\n        # formatting of the span's first run (the Courier-New \"<id>\" run),\n        # then we put the final text back. Assigning the identical text\n        # directly would be a no-op and leave the three runs untouched.\n        $r.Text = \"TEMP_PLACEHOLDER_\" + $id\n        $r.Text = $target\n    }\n}\n"}
